$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" conversion text (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.36 = 49654.22 pesos`n✅ 49654.22 pesos = 12.3 = 969.6 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 80.889
$wsTasas.Range("O10").Value = 4016.48
$wsTasas.Range("N12").Value = 4037.99
$wsTasas.Range("O12").Value = 78.84999999999999
